$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column A (Emocion) gets value 0
$ws.Range("A2").Value = 0

# Row 2, column C (Amplitud) changes from the text "xx.xx" to numeric 0
$ws.Range("C2").Value = 0
